$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before the current row 172 (shifts old rows 172-236 down to 174-238)
$ws.Rows.Item(172).Resize(2).Insert()

# Fill new row 172 - copy boilerplate columns from the row below (now row 174, formerly row 172)
# then set the columns that differ for this new record.
$ws.Cells.Item(172, 1).Value = 7
$ws.Cells.Item(172, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(172, 3).Value = "Ñuble"
$ws.Cells.Item(172, 4).Value = (Get-Date -Year 2021 -Month 9 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(172, 4).NumberFormat = $ws.Cells.Item(174, 4).NumberFormat
$ws.Cells.Item(172, 5).Value = 16
$ws.Cells.Item(172, 6).Value = "Fruta"
$ws.Cells.Item(172, 7).Value = 100102
$ws.Cells.Item(172, 8).Value = "Cítricos"
$ws.Cells.Item(172, 9).Value = 100102005
$ws.Cells.Item(172, 10).Value = "Naranja"
$ws.Cells.Item(172, 11).Value = "Navel Late"
$ws.Cells.Item(172, 12).Value = "Primera"
$ws.Cells.Item(172, 13).Value = 240
$ws.Cells.Item(172, 14).Value = 6000
$ws.Cells.Item(172, 15).Value = 6500
$ws.Cells.Item(172, 16).Value = 6250
$ws.Cells.Item(172, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(172, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(172, 19).Value = 417
$ws.Cells.Item(172, 20).Value = 15

# Fill new row 173
$ws.Cells.Item(173, 1).Value = 7
$ws.Cells.Item(173, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(173, 3).Value = "Ñuble"
$ws.Cells.Item(173, 4).Value = (Get-Date -Year 2021 -Month 9 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(173, 4).NumberFormat = $ws.Cells.Item(174, 4).NumberFormat
$ws.Cells.Item(173, 5).Value = 16
$ws.Cells.Item(173, 6).Value = "Fruta"
$ws.Cells.Item(173, 7).Value = 100102
$ws.Cells.Item(173, 8).Value = "Cítricos"
$ws.Cells.Item(173, 9).Value = 100102005
$ws.Cells.Item(173, 10).Value = "Naranja"
$ws.Cells.Item(173, 11).Value = "Navel Late"
$ws.Cells.Item(173, 12).Value = "Segunda"
$ws.Cells.Item(173, 13).Value = 240
$ws.Cells.Item(173, 14).Value = 5000
$ws.Cells.Item(173, 15).Value = 5500
$ws.Cells.Item(173, 16).Value = 5250
$ws.Cells.Item(173, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(173, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(173, 19).Value = 350
$ws.Cells.Item(173, 20).Value = 15

Write-Host "Final dimension:" $ws.UsedRange.Address()
